$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
# A8 holds rich text "Volume 32   Number  5" -> change trailing "5" to "6"
$cellA8 = $ws.Range("A8")
$lenA8 = $cellA8.Value2.Length
$cellA8.Characters($lenA8, 1).Text = "6"

# C9 holds rich text "Report Covering the Week  1/27/2025  Through  2/2/2025"
# -> "Report Covering the Week  2/3/2025  Through  2/9/2025"
$cellC9 = $ws.Range("C9")
$cellC9.Characters(27, 9).Text = "2/3/2025"
$cellC9.Characters(46, 8).Text = "2/9/2025"

# --- Row 14 ---
$ws.Range("N14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N14").Value = -100

# --- Row 15 ---
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 2

# --- Row 16 ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = -46.153846153846
$ws.Range("L16").Value = -41.666666666666
$ws.Range("M16").Value = -66.666666666666
$ws.Range("N16").Value = -93

# --- Row 17 ---
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -30.769230769230
$ws.Range("I17").Value = 19
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = -17.391304347826
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 171.428571428571
$ws.Range("N17").Value = -20.833333333333

# --- Row 18 ---
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -27.777777777777
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -46.428571428571
$ws.Range("M18").Value = -68.085106382978
$ws.Range("N18").Value = -91.017964071856

# --- Row 19 ---
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -2.127659574468
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 67
$ws.Range("K19").Value = -17.910447761194
$ws.Range("L19").Value = -9.836065573770
$ws.Range("M19").Value = -5.172413793103
$ws.Range("N19").Value = -31.25

# --- Row 20 ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = -20.833333333333
$ws.Range("L20").Value = 72.727272727272
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -90.404040404040

# --- Row 21 ---
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -23.076923076923
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -19.811320754717
$ws.Range("I21").Value = 115
$ws.Range("J21").Value = 149
$ws.Range("K21").Value = -22.818791946308
$ws.Range("L21").Value = -14.179104477611
$ws.Range("M21").Value = -27.215189873417
$ws.Range("N21").Value = -79.930191972076

# --- Row 22 ---
$ws.Range("C22").NumberFormat = 'General'
$ws.Range("C22").Value = "'0"

# --- Row 24 ---
$ws.Range("C24").Value = 34
$ws.Range("E24").Value = -8.108108108108
$ws.Range("F24").Value = 123
$ws.Range("G24").Value = 156
$ws.Range("H24").Value = -21.153846153846
$ws.Range("I24").Value = 171
$ws.Range("J24").Value = 213
$ws.Range("K24").Value = -19.718309859154
$ws.Range("L24").Value = -1.724137931034
$ws.Range("M24").Value = 72.727272727272

# --- Row 25 ---
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 94
$ws.Range("H25").Value = 1.063829787234
$ws.Range("I25").Value = 121
$ws.Range("J25").Value = 124
$ws.Range("K25").Value = -2.419354838709
$ws.Range("L25").Value = 15.238095238095

# --- Row 26 ---
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = -65.384615384615
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = -22.807017543859
$ws.Range("I26").Value = 60
$ws.Range("J26").Value = 90
$ws.Range("K26").Value = -33.333333333333
$ws.Range("L26").Value = 15.384615384615
$ws.Range("M26").Value = 33.333333333333

# --- Row 27 ---
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 5

# --- Row 28 ---
$ws.Range("C28").NumberFormat = 'General'
$ws.Range("C28").Value = "'0"
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -12.5

# --- Row 29 ---
$ws.Range("C29").NumberFormat = 'General'
$ws.Range("C29").Value = "'0"

# --- Row 30 ---
$ws.Range("C30").NumberFormat = 'General'
$ws.Range("C30").Value = "'0"

# --- Row 31 ---
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E31").Value = -100
$ws.Range("G31").NumberFormat = '#,##0'
$ws.Range("G31").Value = 1
$ws.Range("H31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H31").Value = -100
$ws.Range("J31").NumberFormat = '#,##0'
$ws.Range("J31").Value = 1
$ws.Range("K31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K31").Value = -100
